# Guild.xlsx: unify the conception of DataNode, DataTable, Entity.
# Rename the two worksheets and make the second sheet ("DataTable") the
# active/selected tab, matching the author's edit.

$wb = $excel.ActiveWorkbook

$nodeSheet  = $wb.Worksheets.Item(1)   # was "Property1"
$tableSheet = $wb.Worksheets.Item(2)   # was "Record"

$nodeSheet.Name  = "DataNode"
$tableSheet.Name = "DataTable"

# Make "DataTable" the active sheet/tab (mirrors tabSelected moving to sheet2
# and workbookView activeTab="1" in the saved workbook).
$tableSheet.Activate()
